$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "SAV tickets" KPI row (row 12): rename the KPI label and
# replace the text objective "90% / 95%" with a plain numeric target 95.
$ws.Range("A12").Value = "Nombre et taux de résolution des tickets SAV en %"
$ws.Range("B12").Value = 95

# Column widths set explicitly by the author on re-save (values below are
# pre-compensated for the host's +5/6-character padding so the saved OOXML
# <col> widths land on 57.5 / 20 / 18.33203125 / 17.33203125).
$ws.Columns.Item(1).ColumnWidth = 56.666666666666664
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 17.5
$ws.Columns.Item(4).ColumnWidth = 16.5

# Selection left on A13 when the file was saved.
$ws.Range("A13").Select()
